$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Hello" prefix with "World" prefix in column A (rows 2-10)
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    $cell.Value2 = $current -replace '^Hello', 'World'
}

# Update the selected cell / active cell in the sheet view
$ws.Range("C16").Select()
